$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Fix the misspelled survey-response labels that were used throughout the
# answer grid: "Agreee" -> "Agree" and "Strongly agreee" -> "Strongly agree".
# The typo appears in the Likert-scale answer columns (B:J and O) for every
# student row (rows 2-21).
$usedRange = $ws.UsedRange
$rowCount = $usedRange.Rows.Count
$colCount = $usedRange.Columns.Count

for ($r = 1; $r -le $rowCount; $r++) {
    for ($c = 1; $c -le $colCount; $c++) {
        $cell = $usedRange.Cells.Item($r, $c)
        $val = $cell.Value()
        if ($val -eq "Agreee") {
            $cell.Value = "Agree"
        } elseif ($val -eq "Strongly agreee") {
            $cell.Value = "Strongly agree"
        }
    }
}

# Restore the last active selection cell recorded in the saved file
$ws.Range("C3").Select()

$wb.Save()
